$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matches how the source data was stored).
$textCells = @('D5','D6','D9','D10','D11','D12','D16','D17','D20','D21','D22','D23','D24','D25','D27','D29','D30','D31','D32','D33','D35','D36','D38','D39','D40','D41','D42','D43','D45','D46','D47','D48','D49','D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.971.05'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '2.227.67'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  -1.09%  '
$ws.Range('D5').Value = '299.64'
$ws.Range('E5').Value = '  -3.33%  '
$ws.Range('D6').Value = '90.50'
$ws.Range('E6').Value = '  -5.03%  '
$ws.Range('E7').Value = '  -3.69%  '
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D9').Value = '0.491'
$ws.Range('E9').Value = '  -7.08%  '
$ws.Range('D10').Value = '33.06'
$ws.Range('E10').Value = '  -5.86%  '
$ws.Range('D11').Value = '0.0779'
$ws.Range('E11').Value = '  -4.08%  '
$ws.Range('D12').Value = '6.95'
$ws.Range('E12').Value = '  -4.87%  '
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('D14').Value = '2.567.54'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').Value = '2.212.98'
$ws.Range('E15').Value = '  -5.00%  '
$ws.Range('D16').Value = '13.46'
$ws.Range('E16').Value = '  -1.76%  '
$ws.Range('D17').Value = '0.777'
$ws.Range('E17').Value = '  -7.84%  '
$ws.Range('D18').Value = '43.853.12'
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').Value = '0.0₃0903'
$ws.Range('E19').Value = '  -6.57%  '
$ws.Range('D20').Value = '5.94'
$ws.Range('D21').Value = '11.22'
$ws.Range('E21').Value = '  -8.19%  '
$ws.Range('D22').Value = '64.70'
$ws.Range('E22').Value = '  -1.89%  '
$ws.Range('D23').Value = '237.23'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('D24').Value = '2.82'
$ws.Range('E24').Value = '  -6.15%  '
$ws.Range('D25').Value = '0.999'
$ws.Range('E26').Value = '  -6.95%  '
$ws.Range('D27').Value = '38.31'
$ws.Range('E27').Value = '  +1.54%  '
$ws.Range('E28').Value = '  -2.47%  '
$ws.Range('D29').Value = '9.32'
$ws.Range('E29').Value = '  -5.57%  '
$ws.Range('D30').Value = '19.28'
$ws.Range('E30').Value = '  -4.43%  '
$ws.Range('D31').Value = '150.67'
$ws.Range('E31').Value = '  -1.39%  '
$ws.Range('D32').Value = '5.41'
$ws.Range('E32').Value = '  -10.35%  '
$ws.Range('D33').Value = '0.0752'
$ws.Range('E33').Value = '  -6.94%  '
$ws.Range('E34').Value = '  -5.60%  '
$ws.Range('D35').Value = '0.115'
$ws.Range('E35').Value = '  -4.11%  '
$ws.Range('D36').Value = '2.83'
$ws.Range('E36').Value = '  -10.92%  '
$ws.Range('E37').Value = '  -7.57%  '
$ws.Range('D38').Value = '1.70'
$ws.Range('E38').Value = '  -6.42%  '
$ws.Range('D39').Value = '0.0301'
$ws.Range('E39').Value = '  -0.58%  '
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').Value = '3.21'
$ws.Range('E40').Value = '  -7.05%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '3.64'
$ws.Range('E41').Value = '  -5.10%  '
$ws.Range('D42').Value = '13.19'
$ws.Range('E42').Value = '  -8.78%  '
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.87%  '
$ws.Range('D44').Value = '1.839.21'
$ws.Range('E44').Value = '  +4.89%  '
$ws.Range('D45').Value = '1.79'
$ws.Range('E45').Value = '  +12.43%  '
$ws.Range('D46').Value = '0.180'
$ws.Range('E46').Value = '  -7.12%  '
$ws.Range('B47').Value = 'ordi'
$ws.Range('C47').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D47').Value = '67.54'
$ws.Range('E47').Value = '  -4.91%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '14.13'
$ws.Range('E48').Value = '  +9.54%  '
$ws.Range('D49').Value = '73.41'
$ws.Range('E49').Value = '  -9.33%  '
$ws.Range('D50').Value = '94.06'
$ws.Range('E50').Value = '  -6.03%  '
$ws.Range('D51').Value = '2.449.39'
$ws.Range('E51').Value = '  -1.03%  '
